# "alteracoes nos meses para ingles"
#
# The month abbreviations in column A (rows 2:22) of both the DEC and
# FEC sheets are translated from Portuguese to English. jan/mar/jun/jul/nov
# already match between the two languages, so only the other seven
# abbreviations actually change text.

$wb = $excel.ActiveWorkbook

$monthMap = @{
    "fev/23" = "feb/23";
    "abr/23" = "apr/23";
    "mai/23" = "may/23";
    "ago/23" = "aug/23";
    "set/23" = "sep/23";
    "out/23" = "oct/23";
    "dez/23" = "dec/23";
    "fev/24" = "feb/24";
    "abr/24" = "apr/24";
    "mai/24" = "may/24";
    "ago/24" = "aug/24";
    "set/24" = "sep/24";
}

foreach ($sheetName in @("DEC", "FEC")) {
    $ws = $wb.Worksheets.Item($sheetName)
    for ($r = 2; $r -le 22; $r++) {
        $cell = $ws.Cells.Item($r, 1)
        $cur = $cell.Value2
        if ($monthMap.ContainsKey($cur)) {
            $cell.Value = $monthMap[$cur]
        }
    }
}

$dec = $wb.Worksheets.Item("DEC")
$fec = $wb.Worksheets.Item("FEC")

# The trailing blank formatted row below the table moves from DEC down to
# FEC (row 23 disappears from DEC, reappears on FEC).
[void]$dec.Rows.Item(23).Delete()
$fec.Cells.Item(23, 1).Font.Size = 8

# The user ends up on the DEC tab (it becomes the active sheet) with C12
# selected there, while FEC is left with its A1:A22 selection.
[void]$dec.Activate()
[void]$dec.Range("C12").Select()
[void]$fec.Range("A1:A22").Select()
[void]$dec.Activate()
